# The deck's single Design ("Integral" / "Red Violet" colour scheme, stored
# in ppt/theme/theme1.xml) is switched to the stock PowerPoint "Office"
# theme colour scheme (the palette that, before this edit, only lived in the
# otherwise-unused ppt/theme/theme2.xml backing the Notes Master).
#
# Office theme colours (RRGGBB), in clrScheme order:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72
#
# PowerPoint's RGB long is 0xBBGGRR, so the values below are the
# byte-reversed (B*65536 + G*256 + R) form of each hex triplet.

$p = $ppt.ActivePresentation

$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$range = $p.Slides.Range()
$scheme = $range.ThemeColorScheme
for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = $officeColors[$i - 1]
}
